$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.485
$ws.Range("A7").Value = -19.657
$ws.Range("C7").Value = -12.748
$ws.Range("C15").Value = -13.491
$ws.Range("A16").Value = -21.929
$ws.Range("E16").Value = 16.683
$ws.Range("E19").Value = 16.45
$ws.Range("C21").Value = -12.629
$ws.Range("C22").Value = -12.867
$ws.Range("C23").Value = -12.141
$ws.Range("A28").Value = -21.928
$ws.Range("A29").Value = -21.344
$ws.Range("A32").Value = -21.904
$ws.Range("C34").Value = -11.808
$ws.Range("E36").Value = 16.696
$ws.Range("A40").Value = -19.935
$ws.Range("C43").Value = -12.759
$ws.Range("C45").Value = -13.03
$ws.Range("E46").Value = 16.934
$ws.Range("C50").Value = -13.497
$ws.Range("E50").Value = 16.542
$ws.Range("C51").Value = -10.726
$ws.Range("A52").Value = -21.885
$ws.Range("A57").Value = -22.179
$ws.Range("A66").Value = -21.645
$ws.Range("C66").Value = -10.916
$ws.Range("C67").Value = -11.34
$ws.Range("C79").Value = -11.729
$ws.Range("C84").Value = -13.72
$ws.Range("C92").Value = -11.084
$ws.Range("E95").Value = 17.329
$ws.Range("C97").Value = -11.87
$ws.Range("E97").Value = 17.089
$ws.Range("A100").Value = -22.067
